$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-05-31 Friday" "2024-06-01 Saturday"

Replace-Text "623÷8=" "304÷4="
Replace-Text "756÷5=" "844÷4="
Replace-Text "730÷8=" "735÷7="
Replace-Text "397÷8=" "569÷7="
Replace-Text "938÷3=" "760÷4="
Replace-Text "377÷3=" "499÷8="
Replace-Text "964÷8=" "572÷5="
Replace-Text "516÷9=" "749÷3="
Replace-Text "259÷8=" "160÷3="
Replace-Text "711÷5=" "599÷9="
Replace-Text "544÷5=" "773÷9="
Replace-Text "747÷2=" "628÷7="
Replace-Text "534÷8=" "298÷9="
Replace-Text "895÷6=" "708÷2="
Replace-Text "900÷3=" "608÷7="
Replace-Text "920÷5=" "519÷5="
Replace-Text "449÷5=" "860÷5="
Replace-Text "854÷3=" "534÷5="
Replace-Text "463÷9=" "339÷2="
Replace-Text "625÷8=" "480÷7="
Replace-Text "350÷7=" "808÷9="
Replace-Text "156÷3=" "130÷2="
Replace-Text "192÷2=" "875÷6="
Replace-Text "685÷9=" "666÷2="
Replace-Text "360÷9=" "505÷7="

Write-Output "Done"
